# Refined metadata to be additional tab
#
# 1) Update the "time_taken" (column F) timestamps on the existing "data" sheet
# 2) Add a new "metadata" worksheet (after "data") describing the panel query

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update time_taken (F column) timestamps on "data" sheet ---
$ws.Range("F2").Value  = "2021-10-05 14:21:50.958973"
$ws.Range("F3").Value  = "2021-10-05 14:21:50.958981"
$ws.Range("F4").Value  = "2021-10-05 14:21:50.958984"
$ws.Range("F5").Value  = "2021-10-05 14:21:50.958986"
$ws.Range("F6").Value  = "2021-10-05 14:21:50.958989"
$ws.Range("F7").Value  = "2021-10-05 14:21:50.958992"
$ws.Range("F8").Value  = "2021-10-05 14:21:50.958994"
$ws.Range("F9").Value  = "2021-10-05 14:21:50.958997"
$ws.Range("F10").Value = "2021-10-05 14:21:50.958999"
$ws.Range("F11").Value = "2021-10-05 14:21:50.959002"
$ws.Range("F12").Value = "2021-10-05 14:21:50.959004"
$ws.Range("F13").Value = "2021-10-05 14:21:50.959007"
$ws.Range("F14").Value = "2021-10-05 14:21:50.959009"
$ws.Range("F15").Value = "2021-10-05 14:21:50.959012"
$ws.Range("F16").Value = "2021-10-05 14:21:50.959014"
$ws.Range("F17").Value = "2021-10-05 14:21:50.959017"
$ws.Range("F18").Value = "2021-10-05 14:21:50.959019"

# --- 2. Add the new "metadata" worksheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1) - values first, then copy the header style (bold/border)
# used by the "data" sheet's header row so the look matches.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$ws.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row 2
$meta.Range("B2").Value = "Non-CF bronchiectasis"
$meta.Range("C2").Value = 296
$meta.Range("E2").Value = "2020-10-05T14:16:52.566813Z"
$meta.Range("F2").Value = "2021-10-05 14:21:50.955219"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/296/?format=json"

# D2 ("1.26") must stay textual rather than becoming the number 1.26, but
# keep the default (unstyled) cell format - force text, set value, then
# drop back to the Normal style so no numFmt/style id lingers on the cell.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.26"
$meta.Range("D2").Style = "Normal"

# A2 is a plain 0, styled the same way as A2 on the "data" sheet.
$meta.Range("A2").Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("A1").Select()
